$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 120 : LeetCode 3363 - Find the Maximum Number of Fruits Collected
# (no Notes/J entry -> copy a template row that also has no J entry)
# ---------------------------------------------------------------------------
$ws.Range("A118:I118").Copy($ws.Range("A120:I120"))
$ws.Range("J120").Clear()

$ws.Range("A120").Value2 = 3363
$ws.Range("B120").Value2 = "Find the Maximum Number of Fruits Collected"
$ws.Range("C120").Value2 = "#dynamic-programming #rolling-array #matrix "
$ws.Range("D120").Value2 = "hard"
$ws.Range("E120").Value2 = 0
$ws.Range("F120").Value2 = 1
$ws.Range("G120").Value2 = 60
$ws.Range("H120").Value2 = (Get-Date -Year 2025 -Month 8 -Day 7 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I120").Value2 = (Get-Date -Year 2025 -Month 8 -Day 7 -Hour 0 -Minute 0 -Second 0)
$ws.Rows.Item(120).RowHeight = 68

# ---------------------------------------------------------------------------
# Row 121 : LeetCode 808 - Soup Servings
# (has a Notes/J entry -> copy a template row that also has a J entry)
# ---------------------------------------------------------------------------
$ws.Range("A119:J119").Copy($ws.Range("A121:J121"))

$ws.Range("A121").Value2 = 808
$ws.Range("B121").Value2 = "Soup Servings"
$ws.Range("C121").Value2 = "#dynamic-programming #math "
$ws.Range("D121").Value2 = "medium"
$ws.Range("E121").Value2 = 0
$ws.Range("F121").Value2 = 1
$ws.Range("G121").Value2 = 20
$ws.Range("H121").Value2 = (Get-Date -Year 2025 -Month 8 -Day 8 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I121").Value2 = (Get-Date -Year 2025 -Month 8 -Day 8 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J121").Value2 = "???"
$ws.Rows.Item(121).RowHeight = 51

# ---------------------------------------------------------------------------
# Row 122 : LeetCode 231 - Power of Two
# (has a Notes/J entry -> copy a template row that also has a J entry)
# ---------------------------------------------------------------------------
$ws.Range("A119:J119").Copy($ws.Range("A122:J122"))

$ws.Range("A122").Value2 = 231
$ws.Range("B122").Value2 = "Power of Two"
$ws.Range("C122").Value2 = "#math #bit-minipulation "
$ws.Range("D122").Value2 = "easy"
$ws.Range("E122").Value2 = 1
$ws.Range("F122").Value2 = 0
$ws.Range("G122").Value2 = 5
$ws.Range("H122").Value2 = (Get-Date -Year 2025 -Month 8 -Day 9 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I122").Value2 = (Get-Date -Year 2025 -Month 8 -Day 9 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J122").Value2 = "return n & (-n) == n"
$ws.Rows.Item(122).RowHeight = 34

# ---------------------------------------------------------------------------
# Row 123 : LeetCode 1198 - Find Smallest Common Element in All Rows
# (no Notes/J entry -> copy a template row that also has no J entry)
# ---------------------------------------------------------------------------
$ws.Range("A118:I118").Copy($ws.Range("A123:I123"))
$ws.Range("J123").Clear()

$ws.Range("A123").Value2 = 1198
$ws.Range("B123").Value2 = "Find Smallest Common Element in All Rows"
$ws.Range("C123").Value2 = "#hash-table #binary-search #matrix #counting "
$ws.Range("D123").Value2 = "medium"
$ws.Range("E123").Value2 = 1
$ws.Range("F123").Value2 = 0
$ws.Range("G123").Value2 = 10
$ws.Range("H123").Value2 = (Get-Date -Year 2025 -Month 8 -Day 9 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I123").Value2 = (Get-Date -Year 2025 -Month 8 -Day 9 -Hour 0 -Minute 0 -Second 0)
$ws.Rows.Item(123).RowHeight = 51

# ---------------------------------------------------------------------------
# Update the active selection to reflect where the author ended up editing
# ---------------------------------------------------------------------------
$ws.Range("E124").Select()
